$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D4").Value = 3.8
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 6.5
$ws.Range("D7").Value = 9
$ws.Range("D8").Value = 12.3

$ws.Range("E4").Select()
